$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.798.49"
$ws.Range("E2").Value = "  +7.21%  "
$ws.Range("D3").Value = "1.776.77"
$ws.Range("E3").Value = "  +4.07%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.558"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.13%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.50"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.83%  "
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0665"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.39%  "
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").Value = "2.031.37"
$ws.Range("E13").Value = "  +4.11%  "
$ws.Range("D14").Value = "1.775.69"
$ws.Range("E14").Value = "  +4.37%  "
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").Value = "33.750.30"
$ws.Range("E16").Value = "  +7.32%  "
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "251.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "0.0₃0738"
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.45%  "
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.17%  "
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0515"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.32%  "
$ws.Range("E35").Value = "  +7.07%  "
$ws.Range("D36").Value = "1.479.33"
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.63%  "
$ws.Range("E38").Value = "  +3.35%  "
$ws.Range("E39").Value = "  +2.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.885"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.59%  "
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0505"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("E46").Value = "  +3.49%  "
$ws.Range("D47").Value = "1.929.33"
$ws.Range("E47").Value = "  +4.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.29%  "
